$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Sort Items by*") {
        $rng = $p.Range

        # Remove the run containing "distance " (leaving the bookmark tags intact).
        $rng.Find.Execute("distance ", $true, $false, $false, $false, $false,
                           $true, 1, $false, "", 2)

        # Change "or price" to "price".
        $rng2 = $p.Range
        $rng2.Find.Execute("or price", $true, $false, $false, $false, $false,
                            $true, 1, $false, "price", 2)
        break
    }
}
